$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Milestone III checkmarks move around (E column "III" markers) ---
# Clear the "III" marks that were removed
$ws.Range("E9").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("E39").ClearContents()

# Add a new "III" mark
$ws.Range("E54").Value = "III"

# Row 65: mark Milestone III complete (F65 = X), now scoring points
$ws.Range("F65").Value = "X"

# Mark Milestone III complete (X) for the "Effective Use of GIT" and
# "All Graphics API Objects cleaned up in memory" carry-over rows
$ws.Range("E83").Value = "X"
$ws.Range("E84").Value = "X"

# --- New source citation links at the bottom of the sheet ---
$ws.Range("A89").Value = "http://www.rastertek.com/dx11tut36.html"

# A90 picks up a (functionally identical) new formatting style when typed,
# so reproduce that by copying the existing citation row's format first and
# then nudging the font so a new cell style record gets created.
$ws.Range("A88").Copy()
$ws.Range("A90").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A90").Value = "https://www.braynzarsoft.net"
$ws.Range("A90").Font.ThemeColor = 1

# --- Selection moves to A90 ---
$ws.Range("A90").Select()
